$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Move "cm(2)" header from F4 to C4 ---
$ws.Range("F4").Copy($ws.Range("C4"))

# --- 2. Column headers in row 5: reuse style-5 cells already in place ---
# B5 needs style s5 -> clone it from C5 (which currently holds "All(3)" with style 5)
$ws.Range("C5").Copy($ws.Range("B5"))
$ws.Range("B5").Value = "Mean Sepal Length"
# C5 already carries style 5, just change its text in place
$ws.Range("C5").Value = "Mean Sepal Width"

# --- 3. Species labels in column A (style s3), sourced from the old label column (B6:B8) ---
# Do this BEFORE overwriting B6:B8 with numeric means below.
$ws.Range("B6").Copy($ws.Range("A6"))
$ws.Range("A6").Value = "Setosa"

$ws.Range("B7").Copy($ws.Range("A7"))
$ws.Range("A7").Value = "Versicolor"

$ws.Range("B8").Copy($ws.Range("A8"))
$ws.Range("A8").Value = "Virginica"

# --- 4. Mean Sepal Length / Mean Sepal Width values (style s6), pulled from the existing
#         "Sepal Length" / "Sepal Width" Mean rows (10 and 12) so the exact stored
#         double representation is preserved. ---
$ws.Range("D10").Copy($ws.Range("B6"))   # Setosa mean sepal length
$ws.Range("D12").Copy($ws.Range("C6"))   # Setosa mean sepal width

$ws.Range("E10").Copy($ws.Range("B7"))   # Versicolor mean sepal length
$ws.Range("E12").Copy($ws.Range("C7"))   # Versicolor mean sepal width

$ws.Range("F10").Copy($ws.Range("B8"))   # Virginica mean sepal length
$ws.Range("F12").Copy($ws.Range("C8"))   # Virginica mean sepal width

# --- 5. Footnote rows move up from 14-18 to 9-12 (note the "(3: ...)" footnote on the
#         old row 17 is dropped, so row 18 follows directly after row 16). ---
$ws.Range("A14").Copy($ws.Range("A9"))
$ws.Range("A15").Copy($ws.Range("A10"))
$ws.Range("A16").Copy($ws.Range("A11"))
$ws.Range("A18").Copy($ws.Range("A12"))

# --- 6. Clean up everything outside the new A1:C12 extent ---
$ws.Range("D1:F18").Clear()
$ws.Range("B9:C13").Clear()
$ws.Range("A13:C18").Clear()

# --- 7. Rich-text fix: embolden the "(1) sepal dimensions" run within A1 ---
$cell = $ws.Range("A1")
$chars = $cell.Characters(10, 20)
$chars.Font.Bold = $true
